$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update RelivePos column (E2:E4) values from "186,6.89,88" to "186,0,88"
$ws.Range("E2").Value = "186,0,88"
$ws.Range("E3").Value = "186,0,88"
$ws.Range("E4").Value = "186,0,88"

# Update the active selection on the sheet from E4 to F7
$ws.Range("F7").Select()
